$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changed from 45204 to 45205 for every data row (2..305)
$ws.Range("C2:C305").Value2 = 45205

# Row 305 gains an explicit row height (matches Excel's custom-height flag)
$ws.Rows.Item(305).RowHeight = 15

# New row 306 with the added record
$ws.Range("A306").Value2 = "A 47706-2023"
$ws.Range("B306").Value2 = 45203
$ws.Range("C306").Value2 = 45205
$ws.Range("D306").Value2 = "NORRBOTTENS LÄN"
$ws.Range("E306").Value2 = "ÖVERKALIX"
$ws.Range("F306").Value2 = "Sveaskog"
$ws.Range("G306").Value2 = 11.1
$ws.Range("H306").Value2 = 0
$ws.Range("I306").Value2 = 0
$ws.Range("J306").Value2 = 0
$ws.Range("K306").Value2 = 0
$ws.Range("L306").Value2 = 0
$ws.Range("M306").Value2 = 0
$ws.Range("N306").Value2 = 0
$ws.Range("O306").Value2 = 0
$ws.Range("P306").Value2 = 0
$ws.Range("Q306").Value2 = 0
$ws.Range("R306").Value2 = ""

# Formats for the new row matching the rest of the table
$ws.Range("B306:C306").NumberFormat = "YYYY-MM-DD"
$ws.Range("R306").WrapText = $true
